$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (student 1)
$ws.Range("B4").Value = "DTC155D4801030048"
$ws.Range("C4").Value = "Kỹ thuật phần mềm"
$ws.Range("H4").Value = "0916854487"
$ws.Range("H5").Value = "0369548757"
$ws.Range("N4").Value = "085522485"
$ws.Range("N5").Value = "085535125"
$ws.Range("F4").Value = "Nam"
$ws.Range("F5").Value = "Nam"
$ws.Range("U4").Value = "vvc@gmail.com"

# Update the active selection / view position to match the recorded state
$ws.Range("U4").Select()
